$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the four measured threshold values (rows 2-3, columns B/C)
$ws.Range("B2").Value = 5.8
$ws.Range("C2").Value = 10.7
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 9.3

# Move the active selection from F4 to C4
$ws.Range("C4").Select()
